# Add the new "2021年" row (row 5) to the sheet, matching the formatting
# already used by the previous year row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 4's formatting onto row 5 so the new label cell (A5) picks up
# the same bold/bordered/centered style used by the other year labels.
$ws.Range("A4:M4").Copy()
$ws.Range("A5:M5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row label
$ws.Range("A5").Value = "2021年"

# Data values for 2021年
$ws.Range("B5").Value = 62.264
$ws.Range("C5").Value = 55.373
$ws.Range("D5").Value = 37.227
$ws.Range("E5").Value = 43.619
$ws.Range("F5").Value = 56.435
$ws.Range("G5").Value = 55.392
$ws.Range("H5").Value = 56.592
$ws.Range("I5").Value = 54.932
$ws.Range("J5").Value = 43.333
$ws.Range("K5").Value = 55.659
$ws.Range("L5").Value = 82.62
$ws.Range("M5").Value = 26.204
